$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") extending the table.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Give the new headers the same formatting as the existing header row
# (bold, centered/top-aligned, bordered) by copying H1's formats over.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data columns I and J for rows 2-5.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 2

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 7
